$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.747119
$ws.Range("H2").Value = 2.241357
$ws.Range("I2").Value = 0.03096954854571248
$ws.Range("J2").Value = 0.03096954854571248
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 23.59622066666667
$ws.Range("N2").Value = 70.788662
$ws.Range("O2").Value = 0.6996728317814862
$ws.Range("P2").Value = 0.6996728317814862
$ws.Range("Q2").Value = 17.62918478825933
$ws.Range("R2").Value = 158.662663094334
$ws.Range("S2").Value = 0.02166855172997286
$ws.Range("T2").Value = 0.02166855172997286

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.747119
$ws.Range("H3").Value = 2.241357
$ws.Range("I3").Value = 0.03096954854571248
$ws.Range("J3").Value = 0.03096954854571248
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 7.778025666666667
$ws.Range("N3").Value = 23.334077
$ws.Range("O3").Value = 0.2306332577891816
$ws.Range("P3").Value = 0.2306332577891816
$ws.Range("Q3").Value = 5.811110758054333
$ws.Range("R3").Value = 52.299996822489
$ws.Range("S3").Value = 0.007142607873357881
$ws.Range("T3").Value = 0.007142607873357881

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.747119
$ws.Range("H4").Value = 2.241357
$ws.Range("I4").Value = 0.03096954854571248
$ws.Range("J4").Value = 0.03096954854571248
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 2.350402666666667
$ws.Range("N4").Value = 7.051208000000001
$ws.Range("O4").Value = 0.06969391042933218
$ws.Range("P4").Value = 0.06969391042933218
$ws.Range("Q4").Value = 1.756030489917334
$ws.Range("R4").Value = 15.804274409256
$ws.Range("S4").Value = 0.00215838894238174
$ws.Range("T4").Value = 0.00215838894238174

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 19.74619233333334
$ws.Range("H5").Value = 59.23857700000001
$ws.Range("I5").Value = 0.8185184181638298
$ws.Range("J5").Value = 0.8185184181638298
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 23.59622066666667
$ws.Range("N5").Value = 70.788662
$ws.Range("O5").Value = 0.6996728317814862
$ws.Range("P5").Value = 0.6996728317814862
$ws.Range("Q5").Value = 465.935511623775
$ws.Range("R5").Value = 4193.419604613975
$ws.Range("S5").Value = 0.5726950995019895
$ws.Range("T5").Value = 0.5726950995019895

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 19.74619233333334
$ws.Range("H6").Value = 59.23857700000001
$ws.Range("I6").Value = 0.8185184181638298
$ws.Range("J6").Value = 0.8185184181638298
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 7.778025666666667
$ws.Range("N6").Value = 23.334077
$ws.Range("O6").Value = 0.2306332577891816
$ws.Range("P6").Value = 0.2306332577891816
$ws.Range("Q6").Value = 153.5863907876032
$ws.Range("R6").Value = 1382.277517088429
$ws.Range("S6").Value = 0.1887775693415717
$ws.Range("T6").Value = 0.1887775693415717

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 19.74619233333334
$ws.Range("H7").Value = 59.23857700000001
$ws.Range("I7").Value = 0.8185184181638298
$ws.Range("J7").Value = 0.8185184181638298
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.350402666666667
$ws.Range("N7").Value = 7.051208000000001
$ws.Range("O7").Value = 0.06969391042933218
$ws.Range("P7").Value = 0.06969391042933218
$ws.Range("Q7").Value = 46.41150311677957
$ws.Range("R7").Value = 417.7035280510161
$ws.Range("S7").Value = 0.05704574932026861
$ws.Range("T7").Value = 0.05704574932026861

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3.630999
$ws.Range("H8").Value = 10.892997
$ws.Range("I8").Value = 0.1505120332904577
$ws.Range("J8").Value = 0.1505120332904577
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 23.59622066666667
$ws.Range("N8").Value = 70.788662
$ws.Range("O8").Value = 0.6996728317814862
$ws.Range("P8").Value = 0.6996728317814862
$ws.Range("Q8").Value = 85.67785364444599
$ws.Range("R8").Value = 771.1006828000139
$ws.Range("S8").Value = 0.1053091805495238
$ws.Range("T8").Value = 0.1053091805495239

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3.630999
$ws.Range("H9").Value = 10.892997
$ws.Range("I9").Value = 0.1505120332904577
$ws.Range("J9").Value = 0.1505120332904577
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 7.778025666666667
$ws.Range("N9").Value = 23.334077
$ws.Range("O9").Value = 0.2306332577891816
$ws.Range("P9").Value = 0.2306332577891816
$ws.Range("Q9").Value = 28.242003417641
$ws.Range("R9").Value = 254.178030758769
$ws.Range("S9").Value = 0.03471308057425201
$ws.Range("T9").Value = 0.03471308057425202

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.630999
$ws.Range("H10").Value = 10.892997
$ws.Range("I10").Value = 0.1505120332904577
$ws.Range("J10").Value = 0.1505120332904577
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2.350402666666667
$ws.Range("N10").Value = 7.051208000000001
$ws.Range("O10").Value = 0.06969391042933218
$ws.Range("P10").Value = 0.06969391042933218
$ws.Range("Q10").Value = 8.534309732264001
$ws.Range("R10").Value = 76.808787590376
$ws.Range("S10").Value = 0.01048977216668182
$ws.Range("T10").Value = 0.01048977216668182
